$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 8.956950000000001
$ws.Cells.Item(2, 8).Value = 26.87085
$ws.Cells.Item(2, 9).Value = 0.3465211830970586
$ws.Cells.Item(2, 10).Value = 0.3465211830970586
$ws.Cells.Item(2, 13).Value = 0.1356863333333333
$ws.Cells.Item(2, 14).Value = 0.407059
$ws.Cells.Item(2, 15).Value = 0.004454204096299941
$ws.Cells.Item(2, 16).Value = 0.004454204096299941
$ws.Cells.Item(2, 17).Value = 1.21533570335
$ws.Cells.Item(2, 18).Value = 10.93802133015
$ws.Cells.Item(2, 19).Value = 0.00154347607320562
$ws.Cells.Item(2, 20).Value = 0.00154347607320562
$ws.Cells.Item(3, 7).Value = 8.956950000000001
$ws.Cells.Item(3, 8).Value = 26.87085
$ws.Cells.Item(3, 9).Value = 0.3465211830970586
$ws.Cells.Item(3, 10).Value = 0.3465211830970586
$ws.Cells.Item(3, 15).Value = 0.564494940478519
$ws.Cells.Item(3, 16).Value = 0.5644949404785189
$ws.Cells.Item(3, 17).Value = 154.02321956775
$ws.Cells.Item(3, 18).Value = 1386.20897610975
$ws.Cells.Item(3, 19).Value = 0.1956094546269201
$ws.Cells.Item(3, 20).Value = 0.19560945462692
$ws.Cells.Item(4, 7).Value = 8.956950000000001
$ws.Cells.Item(4, 8).Value = 26.87085
$ws.Cells.Item(4, 9).Value = 0.3465211830970586
$ws.Cells.Item(4, 10).Value = 0.3465211830970586
$ws.Cells.Item(4, 13).Value = 13.13090033333333
$ws.Cells.Item(4, 14).Value = 39.392701
$ws.Cells.Item(4, 15).Value = 0.4310508554251812
$ws.Cells.Item(4, 16).Value = 0.4310508554251811
$ws.Cells.Item(4, 17).Value = 117.61281774065
$ws.Cells.Item(4, 18).Value = 1058.51535966585
$ws.Cells.Item(4, 19).Value = 0.1493682523969329
$ws.Cells.Item(4, 20).Value = 0.1493682523969329
$ws.Cells.Item(5, 9).Value = 0.2466462208011621
$ws.Cells.Item(5, 10).Value = 0.2466462208011621
$ws.Cells.Item(5, 13).Value = 0.1356863333333333
$ws.Cells.Item(5, 14).Value = 0.407059
$ws.Cells.Item(5, 15).Value = 0.004454204096299941
$ws.Cells.Item(5, 16).Value = 0.004454204096299941
$ws.Cells.Item(5, 17).Value = 0.8650494482238889
$ws.Cells.Item(5, 18).Value = 7.785445034015
$ws.Cells.Item(5, 19).Value = 0.001098612607029436
$ws.Cells.Item(5, 20).Value = 0.001098612607029436
$ws.Cells.Item(6, 9).Value = 0.2466462208011621
$ws.Cells.Item(6, 10).Value = 0.2466462208011621
$ws.Cells.Item(6, 15).Value = 0.564494940478519
$ws.Cells.Item(6, 16).Value = 0.5644949404785189
$ws.Cells.Item(6, 19).Value = 0.1392305437304036
$ws.Cells.Item(6, 20).Value = 0.1392305437304036
$ws.Cells.Item(7, 9).Value = 0.2466462208011621
$ws.Cells.Item(7, 10).Value = 0.2466462208011621
$ws.Cells.Item(7, 13).Value = 13.13090033333333
$ws.Cells.Item(7, 14).Value = 39.392701
$ws.Cells.Item(7, 15).Value = 0.4310508554251812
$ws.Cells.Item(7, 16).Value = 0.4310508554251811
$ws.Cells.Item(7, 17).Value = 83.7142386339539
$ws.Cells.Item(7, 18).Value = 753.428147705585
$ws.Cells.Item(7, 19).Value = 0.106317064463729
$ws.Cells.Item(7, 20).Value = 0.106317064463729
$ws.Cells.Item(8, 7).Value = 7.505276333333332
$ws.Cells.Item(8, 8).Value = 22.515829
$ws.Cells.Item(8, 9).Value = 0.2903596910217228
$ws.Cells.Item(8, 10).Value = 0.2903596910217228
$ws.Cells.Item(8, 13).Value = 0.1356863333333333
$ws.Cells.Item(8, 14).Value = 0.407059
$ws.Cells.Item(8, 15).Value = 0.004454204096299941
$ws.Cells.Item(8, 16).Value = 0.004454204096299941
$ws.Cells.Item(8, 17).Value = 1.018363426323444
$ws.Cells.Item(8, 18).Value = 9.165270836910999
$ws.Cells.Item(8, 19).Value = 0.001293321325149343
$ws.Cells.Item(8, 20).Value = 0.001293321325149343
$ws.Cells.Item(9, 7).Value = 7.505276333333332
$ws.Cells.Item(9, 8).Value = 22.515829
$ws.Cells.Item(9, 9).Value = 0.2903596910217228
$ws.Cells.Item(9, 10).Value = 0.2903596910217228
$ws.Cells.Item(9, 15).Value = 0.564494940478519
$ws.Cells.Item(9, 16).Value = 0.5644949404785189
$ws.Cells.Item(9, 17).Value = 129.0603190378016
$ws.Cells.Item(9, 18).Value = 1161.542871340215
$ws.Cells.Item(9, 19).Value = 0.1639065765006686
$ws.Cells.Item(9, 20).Value = 0.1639065765006685
$ws.Cells.Item(10, 7).Value = 7.505276333333332
$ws.Cells.Item(10, 8).Value = 22.515829
$ws.Cells.Item(10, 9).Value = 0.2903596910217228
$ws.Cells.Item(10, 10).Value = 0.2903596910217228
$ws.Cells.Item(10, 13).Value = 13.13090033333333
$ws.Cells.Item(10, 14).Value = 39.392701
$ws.Cells.Item(10, 15).Value = 0.4310508554251812
$ws.Cells.Item(10, 16).Value = 0.4310508554251811
$ws.Cells.Item(10, 17).Value = 98.55103550712543
$ws.Cells.Item(10, 18).Value = 886.9593195641289
$ws.Cells.Item(10, 19).Value = 0.1251597931959049
$ws.Cells.Item(10, 20).Value = 0.1251597931959049
$ws.Cells.Item(11, 7).Value = 3.010615333333333
$ws.Cells.Item(11, 8).Value = 9.031846
$ws.Cells.Item(11, 9).Value = 0.1164729050800565
$ws.Cells.Item(11, 10).Value = 0.1164729050800565
$ws.Cells.Item(11, 13).Value = 0.1356863333333333
$ws.Cells.Item(11, 14).Value = 0.407059
$ws.Cells.Item(11, 15).Value = 0.004454204096299941
$ws.Cells.Item(11, 16).Value = 0.004454204096299941
$ws.Cells.Item(11, 17).Value = 0.4084993556571111
$ws.Cells.Item(11, 18).Value = 3.676494200914
$ws.Cells.Item(11, 19).Value = 0.0005187940909155419
$ws.Cells.Item(11, 20).Value = 0.0005187940909155419
$ws.Cells.Item(12, 7).Value = 3.010615333333333
$ws.Cells.Item(12, 8).Value = 9.031846
$ws.Cells.Item(12, 9).Value = 0.1164729050800565
$ws.Cells.Item(12, 10).Value = 0.1164729050800565
$ws.Cells.Item(12, 15).Value = 0.564494940478519
$ws.Cells.Item(12, 16).Value = 0.5644949404785189
$ws.Cells.Item(12, 17).Value = 51.77037568815666
$ws.Cells.Item(12, 18).Value = 465.93338119341
$ws.Cells.Item(12, 19).Value = 0.06574836562052669
$ws.Cells.Item(12, 20).Value = 0.06574836562052667
$ws.Cells.Item(13, 7).Value = 3.010615333333333
$ws.Cells.Item(13, 8).Value = 9.031846
$ws.Cells.Item(13, 9).Value = 0.1164729050800565
$ws.Cells.Item(13, 10).Value = 0.1164729050800565
$ws.Cells.Item(13, 13).Value = 13.13090033333333
$ws.Cells.Item(13, 14).Value = 39.392701
$ws.Cells.Item(13, 15).Value = 0.4310508554251812
$ws.Cells.Item(13, 16).Value = 0.4310508554251811
$ws.Cells.Item(13, 17).Value = 39.53208988400512
$ws.Cells.Item(13, 18).Value = 355.788808956046
$ws.Cells.Item(13, 19).Value = 0.05020574536861429
$ws.Cells.Item(13, 20).Value = 0.05020574536861427
